$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.689.83"
Set-TextValue "E2" "  -3.11%  "
Set-TextValue "D3" "2.095.98"
Set-TextValue "E3" "  -1.33%  "
Set-TextValue "D4" "1.008"
Set-TextValue "E4" "  -0.48%  "
Set-TextValue "D5" "343.21"
Set-TextValue "E5" "  -2.82%  "
Set-TextValue "E6" "  -0.45%  "
Set-TextValue "D7" "0.5146"
Set-TextValue "E7" "  -2.47%  "
Set-TextValue "D8" "0.4400"
Set-TextValue "E8" "  -3.14%  "
Set-TextValue "D9" "53.07"
Set-TextValue "E9" "  -1.66%  "
Set-TextValue "D10" "0.09212"
Set-TextValue "E10" "  +1.34%  "
Set-TextValue "D11" "1.170"
Set-TextValue "E11" "  -0.95%  "
Set-TextValue "D12" "24.88"
Set-TextValue "D13" "2.104.17"
Set-TextValue "E13" "  -1.87%  "
Set-TextValue "D14" "6.756"
Set-TextValue "E14" "  -1.37%  "
Set-TextValue "D15" "8.181"
Set-TextValue "E15" "  +1.08%  "
Set-TextValue "D16" "99.29"
Set-TextValue "E16" "  -3.13%  "
Set-TextValue "D17" "0.00001150"
Set-TextValue "E17" "  -2.17%  "
Set-TextValue "D18" "1.008"
Set-TextValue "E18" "  -0.48%  "
Set-TextValue "D19" "20.71"
Set-TextValue "E19" "  +6.47%  "
Set-TextValue "D20" "0.06632"
Set-TextValue "E20" "  -1.16%  "
Set-TextValue "E21" "  -0.48%  "
Set-TextValue "D22" "6.188"
Set-TextValue "E22" "  -2.44%  "
Set-TextValue "D23" "29.742.65"
Set-TextValue "E23" "  -3.15%  "
Set-TextValue "D24" "12.55"
Set-TextValue "E24" "  -2.60%  "
Set-TextValue "D25" "2.298"
Set-TextValue "E25" "  -3.90%  "
Set-TextValue "D26" "2.348.67"
Set-TextValue "E26" "  -1.42%  "
Set-TextValue "D27" "21.83"
Set-TextValue "E27" "  -2.95%  "
Set-TextValue "D28" "162.80"
Set-TextValue "E28" "  -1.20%  "
Set-TextValue "D29" "2.521"
Set-TextValue "E29" "  -1.77%  "
Set-TextValue "D30" "132.53"
Set-TextValue "E30" "  -2.73%  "
Set-TextValue "D31" "1.132"
Set-TextValue "E31" "  -5.37%  "
Set-TextValue "E32" "  -3.08%  "
Set-TextValue "D33" "1.637"
Set-TextValue "E33" "  -1.14%  "
Set-TextValue "D34" "6.152"
Set-TextValue "E34" "  -3.30%  "
Set-TextValue "D35" "3.963"
Set-TextValue "E35" "  -1.43%  "
Set-TextValue "D36" "6.053"
Set-TextValue "E36" "  -1.63%  "
Set-TextValue "D37" "10.23"
Set-TextValue "E37" "  -0.86%  "
Set-TextValue "D38" "0.02561"
Set-TextValue "E38" "  -3.54%  "
Set-TextValue "D39" "0.06700"
Set-TextValue "E39" "  -2.76%  "
Set-TextValue "D40" "12.39"
Set-TextValue "E40" "  -1.29%  "
Set-TextValue "D41" "0.6852"
Set-TextValue "E41" "  -1.08%  "
Set-TextValue "D42" "0.2215"
Set-TextValue "E42" "  -4.54%  "
Set-TextValue "D43" "1.292"
Set-TextValue "E43" "  +0.98%  "
Set-TextValue "D44" "0.6624"
Set-TextValue "E44" "  +2.63%  "
Set-TextValue "D45" "14.11"
Set-TextValue "E45" "  -4.66%  "
Set-TextValue "E46" "  -2.04%  "
Set-TextValue "D47" "3.611"
Set-TextValue "E47" "  -4.29%  "
Set-TextValue "E48" "  -3.14%  "
Set-TextValue "D49" "0.00000000337"
Set-TextValue "E49" "  -8.45%  "
Set-TextValue "D50" "81.90"
Set-TextValue "E50" "  -1.36%  "
Set-TextValue "B51" "WOONetwork"
Set-TextValue "C51" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D51" "0.3259"
Set-TextValue "E51" "  -1.86%  "
